# Add two new trailing columns ("ParcelSize", "Weight") to the header row,
# size the new columns to fit their header text (mirrors the existing
# "best fit" columns A:J), and move the active selection to the next free
# cell past the new data, the way Excel leaves it after you type the last
# header and hit Enter/Tab off the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "ParcelSize"
$ws.Range("L1").Value = "Weight"

# Best-fit the two new columns to their header text, same as columns A:J.
$ws.Columns.Item(11).ColumnWidth = 9.140625
$ws.Columns.Item(12).ColumnWidth = 6.7

# Leave the selection where Excel would land after finishing data entry on
# the header row (one cell past the last populated column).
$ws.Range("M1").Select() | Out-Null
